# Generate Report for handback
# Update "Correspond Handoff Datetime" (col D) and "Correspond Handback DateTime" (col G)
# for the 5fe91086-... handback row (row 3) on both language sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-08 14:17:44"
$wsZhCn.Range("G3").Value = "2016-01-08 14:18:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-08 14:17:59"
$wsDeDe.Range("G3").Value = "2016-01-08 14:19:09"
